# Fruta / hortaliza, semanal
# Insert this week's new observation (2023-10-06, serial 45205) as a new
# row right after the current last-but-one block entry (old row 34), which
# pushes the existing rows 35-51 down to 36-52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(35).Insert()

$ws.Cells.Item(35, 1).Value  = 8
$ws.Cells.Item(35, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(35, 3).Value  = "Coquimbo"
$ws.Cells.Item(35, 4).Value  = 45205
$ws.Cells.Item(35, 5).Value  = 4
$ws.Cells.Item(35, 6).Value  = "Fruta"
$ws.Cells.Item(35, 7).Value  = 100101
$ws.Cells.Item(35, 8).Value  = "Berries"
$ws.Cells.Item(35, 9).Value  = 100101001
$ws.Cells.Item(35, 10).Value = "Arándano (blue)"
$ws.Cells.Item(35, 11).Value = "Sin especificar"
$ws.Cells.Item(35, 12).Value = "Primera"
$ws.Cells.Item(35, 13).Value = 200
$ws.Cells.Item(35, 14).Value = 12000
$ws.Cells.Item(35, 15).Value = 13000
$ws.Cells.Item(35, 16).Value = 12500
$ws.Cells.Item(35, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(35, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(35, 19).Value = 6250
$ws.Cells.Item(35, 20).Value = 2
